$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Split the old "Abstract" style into "Abstract Title" (new) and a
#    revised "Abstract" body style (Pandoc 3 reference.docx behaviour).
# ------------------------------------------------------------------

# 1a. Tweak the existing "Abstract" style first (before we add the new
#     "Abstract Title" style) so the later lookups by name stay simple.
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# 1b. Add the new "Abstract Title" style, based on Normal, followed by
#     the (now revised) "Abstract" style.
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = $d.Styles("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles("Abstract")
$abstractTitle.QuickStyle = $true
$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060
$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10

# ------------------------------------------------------------------
# 2. Add a new "Footnote Block Text" style (block-quote text inside
#    footnotes), based on "Footnote Text".
# ------------------------------------------------------------------
$footnoteBlockText = $d.Styles.Add("Footnote Block Text", 1)
$footnoteBlockText.BaseStyle = $d.Styles("Footnote Text")
$footnoteBlockText.NextParagraphStyle = $d.Styles("Footnote Text")
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true
$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24

# ------------------------------------------------------------------
# 3. Fill in the missing character formatting for the "ImportTok" and
#    "BuiltInTok" syntax-highlighting styles.
# ------------------------------------------------------------------
$importTok = $d.Styles("ImportTok")
$importTok.Font.Bold = $true
$importTok.Font.Color = 32768

$builtInTok = $d.Styles("BuiltInTok")
$builtInTok.Font.Color = 32768
